# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> bound to the (only) slide master, currently the
#                            "Integral" / "Red Violet" palette actually seen
#                            on the slides.
#   ppt/theme/theme2.xml -> bound only to the notes master, currently the
#                            stock "Office Theme" palette (Google Slides
#                            export leftover, never shown on a slide).
#
# The authored edit swaps the two themes' contents, so the slide master
# (and therefore every slide) now renders with the plain "Office" color
# palette instead of "Integral". Apply that by rewriting the 12 theme
# colors on the slide master's theme (theme1.xml) to the Office palette,
# via the slide's ThemeColorScheme, which writes straight through to the
# shared <a:clrScheme> used by every slide/layout.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order is Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink, FollowedHyperlink
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
